$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price-report record needs to be inserted into the weekly series.
# Insert a new row above row 40; this shifts the existing rows 40:143
# down to 41:144 (and grows the used range to A1:R144), matching the
# "Fruta / hortaliza, semanal" update.
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new observation.
$ws.Cells.Item(40, 1).Value  = 4
$ws.Cells.Item(40, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(40, 3).Value  = "Los Lagos"
$ws.Cells.Item(40, 4).Value  = 45114
$ws.Cells.Item(40, 5).Value  = 10
$ws.Cells.Item(40, 6).Value  = 100112026
$ws.Cells.Item(40, 7).Value  = "Haba"
$ws.Cells.Item(40, 8).Value  = "Sin especificar"
$ws.Cells.Item(40, 9).Value  = "Primera"
$ws.Cells.Item(40, 10).Value = 80
$ws.Cells.Item(40, 11).Value = 21000
$ws.Cells.Item(40, 12).Value = 21000
$ws.Cells.Item(40, 13).Value = 21000
$ws.Cells.Item(40, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(40, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(40, 16).Value = 840
$ws.Cells.Item(40, 17).Value = 25
$ws.Cells.Item(40, 18).Value = "Hortaliza"
